$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Header row
$ws.Range("A1").Value = "age"
$ws.Range("B1").Value = "insurance"

# Data rows 2-9: age values and yes/no insurance values
$ages = @(18, 12, 15, 56, 35, 12, 70, 45)
$insurance = @("yes", "no", "no", "yes", "yes", "no", "yes", "yes")

for ($i = 0; $i -lt $ages.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $ages[$i]
    $ws.Cells.Item($row, 2).Value = $insurance[$i]
}

# Rows 10 and 11: clear out column A and B values, leaving only styled empty cells
$ws.Range("A10:B11").ClearContents()

# Update selection to B10
$ws.Range("B10").Select() | Out-Null
